$d = $word.ActiveDocument

$replacements = @(
    @{old = "55×36=1980"; new = "56×97=5432"},
    @{old = "66×73=4818"; new = "94×36=3384"},
    @{old = "14×73=1022"; new = "20×28=560"},
    @{old = "22×58=1276"; new = "99×81=8019"},
    @{old = "20×17=340";  new = "32×86=2752"},
    @{old = "95×69=6555"; new = "12×64=768"},
    @{old = "86×16=1376"; new = "97×47=4559"},
    @{old = "23×53=1219"; new = "23×61=1403"},
    @{old = "47×72=3384"; new = "51×60=3060"},
    @{old = "18×68=1224"; new = "57×98=5586"},
    @{old = "44×76=3344"; new = "38×53=2014"},
    @{old = "88×66=5808"; new = "85×50=4250"},
    @{old = "35×46=1610"; new = "66×71=4686"},
    @{old = "20×39=780";  new = "61×39=2379"},
    @{old = "66×18=1188"; new = "50×79=3950"},
    @{old = "35×91=3185"; new = "46×95=4370"},
    @{old = "41×66=2706"; new = "56×97=5432"},
    @{old = "53×34=1802"; new = "70×24=1680"},
    @{old = "79×78=6162"; new = "85×72=6120"},
    @{old = "19×87=1653"; new = "12×48=576"},
    @{old = "85×92=7820"; new = "17×21=357"},
    @{old = "23×68=1564"; new = "58×72=4176"},
    @{old = "72×72=5184"; new = "82×56=4592"},
    @{old = "76×75=5700"; new = "44×98=4312"},
    @{old = "68×39=2652"; new = "48×67=3216"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
